$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rakib & Mahfuz: 19th (T column) meal = 1.5 each (jelaSomity)
$ws.Range("T3").Value = 1.5
$ws.Range("T4").Value = 1.5

# Update the selection/view to reflect where the user left off editing
$ws.Range("T5").Select()
